$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---- Row 2 (Sourabh Awasthi) ----
$ws1.Range("D2").Value = "sourabh.awasthi@capgemini.com"
$ws1.Range("J2").Value = '#024fL9?"Ynx'

# ---- Row 3 (Sandipan Deb) ----
$ws1.Range("D3").Value = "sandipan.deb@capgemini.com"
$ws1.Range("J3").Value = 'f}Fd/th3Y6(2'

# ---- Sheet2 new rows (username / Company lookups) ----
$ws2.Range("C2").Value = "sandipan.deb"
$ws2.Range("I2").Value = "Capgemini"
$ws2.Range("C3").Value = "sandipan.deb"
$ws2.Range("I3").Value = "Capgemini"

# ---- Row 4 (Biswaji Deb) ----
$ws1.Range("D4").Value = "biswaji.deb@capgemini.com"
$ws1.Range("J4").Value = 'OSgRR5AtJg/['

# ---- Row 5 (Dhiraj Kajari) ----
$ws1.Range("D5").Value = "dhiraj.kajari@capgemini.com"
$ws1.Range("J5").Value = '$02[)"tS#!Cf'

# ---- Row 6 (Manoj Kumar B S) ----
$ws1.Range("D6").Value = "manoj-kumar.b.s@capgemini.com"
$ws1.Range("J6").Value = '&+vc]#h((WpF'

# ---- Row 7 (Mayur Bhorkar) ----
$ws1.Range("D7").Value = "mayur.bhorkar@capgemini.com"
$ws1.Range("J7").Value = 'L!smG}o?gdA?'

# ---- Shared formulas across rows 2:7 ----
$ws1.Range("A2:A7").Formula = "=PROPER(IFERROR(LEFT(C2,FIND(CHAR(46),C2)-1),C2))"
$ws1.Range("B2:B7").Formula = '=IFERROR(PROPER(RIGHT(C2,LEN(C2)-FIND("@",SUBSTITUTE(C2,".","@",((LEN(C2)-LEN(SUBSTITUTE(C2,".","")))/LEN("\")))))), "Unknown")'
$ws1.Range("C2:C7").Formula = "=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D2,FIND(CHAR(64),D2)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))"
$ws1.Range("E2:E7").Formula = "=LEFT(H2,FIND(CHAR(46),H2)-1)"
$ws1.Range("F2:F7").Formula = '=CONCATENATE("ITPartner\",I2)'
$ws1.Range("H2:H7").Formula = "=RIGHT(D2,LEN(D2)-FIND(CHAR(64),D2))"
$ws1.Range("I2:I7").Formula = "=PROPER(E2)"
$ws1.Range("P2:P7").Formula = "=COUNTIF(D:D,D2)"

# ---- Numeric / boolean columns ----
$ws1.Range("K2:K7").Value = 80
$ws1.Range("M2:M7").Value = $true

$excel.Calculate()
